$d = $word.ActiveDocument

# Locate the end of the "jabatanPimpinan" placeholder text.
$jpRange = $d.Content
[void]$jpRange.Find.Execute("jabatanPimpinan")
$jpEnd = $jpRange.End

# Locate the end of the "satker" placeholder text that follows it
# (" ${satker}" comes right after "${jabatanPimpinan}").
$skRange = $d.Content
[void]$skRange.Find.Execute("satker")
$skEnd = $skRange.End

# Remove everything between them: the stray "}" that closed
# "jabatanPimpinan", the space, and the "${satker" opener - leaving just
# "${jabatanPimpinan}" followed by the original closing "}" run.
$toDelete = $d.Range($jpEnd, $skEnd)
$toDelete.Delete()
